$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 6
$ws.Range("A6").Value = 9881
$ws.Range("B6").Value = 10000
$ws.Range("C6").Value = 81.47
$ws.Range("D6").Value = 80.5
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -1.19
$ws.Range("G6").Value = 42607.884201388886
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"
$ws.Range("H6").Value = $false

# Add row 7
$ws.Range("A7").Value = 9835.5499999999993
$ws.Range("B7").Value = 9881
$ws.Range("C7").Value = 80.45
$ws.Range("D7").Value = 80.08
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = -0.46
$ws.Range("G7").Value = 42608.616307870368
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"
$ws.Range("H7").Value = $false
